$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# Version: 0.4.0 -> 0.7.0
$wsMeta.Range("B3").Value = "0.7.0"

# Date: 2023-06-30T18:27:12-04:00 -> 2023-09-13T17:11:14-03:00
$wsMeta.Range("B8").Value = "2023-09-13T17:11:14-03:00"

# Description: drop trailing space (same text reused on Elements!M1)
$description = "Extensión para dato es de autoidentificación, por lo tanto, es una información que la persona entrega y el modelo de atención debe garantizar las condiciones y los mecanismos de privacidad y confidencialidad a través de un protocolo de aplicación local"
$wsMeta.Range("B11").Value = $description
$wsElem.Range("M1").Value = $description

# Context: element:Element -> element:Patient
$wsMeta.Range("B20").Value = "element:Patient"

# Base row for the Extension root element now also documents the ele-1/ext-1 invariants
$wsElem.Range("AJ1").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
